# Updates cryptos list values (prices and volume %) to refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.660.83'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '''1.845.81'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('D4').Value = '''1.011'
$ws.Range('E4').Value = '  -2.79%  '
$ws.Range('D5').Value = '''318.59'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').Value = '''1.009'
$ws.Range('E6').Value = '  -2.66%  '
$ws.Range('D7').Value = '''0.4307'
$ws.Range('E7').Value = '  -2.69%  '
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').Value = '''0.07344'
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('D10').Value = '''0.8803'
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('D11').Value = '''21.55'
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').Value = '''1.849.84'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').Value = '''6.728'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('D14').Value = '''5.451'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = '''0.07094'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '''87.65'
$ws.Range('E16').Value = '  +4.53%  '
$ws.Range('D17').Value = '''1.013'
$ws.Range('E17').Value = '  -2.84%  '
$ws.Range('D18').Value = '''0.000008974'
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('D19').Value = '''1.009'
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').Value = '''15.44'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').Value = '''27.663.41'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = '''5.252'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('E23').Value = '  -1.89%  '
$ws.Range('D24').Value = '''2.079.34'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('D25').Value = '''2.037'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('D26').Value = '''155.35'
$ws.Range('E26').Value = '  -2.07%  '
$ws.Range('D27').Value = '''18.53'
$ws.Range('E27').Value = '  -1.90%  '
$ws.Range('D28').Value = '''2.140'
$ws.Range('E28').Value = '  +7.73%  '
$ws.Range('D29').Value = '''5.384'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').Value = '''120.37'
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').Value = '''0.08908'
$ws.Range('E31').Value = '  -1.91%  '
$ws.Range('D32').Value = '''1.227'
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('D33').Value = '''0.7799'
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('D35').Value = '''2.906'
$ws.Range('E35').Value = '  -6.74%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').Value = '''1.010'
$ws.Range('E36').Value = '  -2.86%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '''1.138'
$ws.Range('E37').Value = '  -1.89%  '
$ws.Range('D38').Value = '''0.05334'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('D39').Value = '''0.01968'
$ws.Range('E39').Value = '  -1.38%  '
$ws.Range('D40').Value = '''7.228'
$ws.Range('E40').Value = '  +4.23%  '
$ws.Range('D41').Value = '''2.865'
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('D42').Value = '''0.5160'
$ws.Range('D43').Value = '''0.1677'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('D44').Value = '''8.925'
$ws.Range('E44').Value = '  +2.73%  '
$ws.Range('D45').Value = '''110.25'
$ws.Range('E45').Value = '  +0.46%  '
$ws.Range('D46').Value = '''10.68'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').Value = '''0.4728'
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('D48').Value = '''0.06497'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D49').Value = '''1.698'
$ws.Range('E49').Value = '  -1.75%  '
$ws.Range('D50').Value = '''1.009'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('D51').Value = '''1.893'
$ws.Range('E51').Value = '  -1.05%  '
